# "ЧМ Pract 3 final fixed"
# Update the manually-entered "C++ solution" reference values (column G,
# rows 6-9) so they line up with Excel's own MMULT results, and leave the
# selection/zoom the way the author left the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- update the "Решение C++:" comparison values (G6:G9) -------------------
# These cells hold their numbers as literal text (not numeric values), so
# the target cell has to be coerced to Text before the assignment, otherwise
# Excel auto-converts the numeric-looking string back into a number. After
# writing the text we restore the original (bordered, General) cell format
# by pasting formats from a neighbouring cell that already carries it, so
# the visual style stays exactly as it was.
$newCppValues = @{
    "G6" = "0.663415"
    "G7" = "0.628239"
    "G8" = "1.65558"
    "G9" = "0.553527"
}

foreach ($addr in $newCppValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $newCppValues[$addr]
}

$ws.Range("F7").Copy() | Out-Null
$ws.Range("G6:G9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- restore the view state (zoom + active selection) ----------------------
$ws.Activate() | Out-Null
$ws.Range("F14").Select() | Out-Null
$excel.ActiveWindow.Zoom = 125
